$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).Value = "68.822.10"
$ws.Cells.Item(2, 5).Value = "  +1.12%  "
$ws.Cells.Item(3, 4).Value = "3.711.75"
$ws.Cells.Item(3, 5).Value = "  +0.05%  "
$ws.Cells.Item(4, 5).Value = "  +0.04%  "
$ws.Cells.Item(5, 4).Value = "616.62"
$ws.Cells.Item(5, 5).Value = "  +5.85%  "
$ws.Cells.Item(6, 4).Value = "188.22"
$ws.Cells.Item(6, 5).Value = "  +5.82%  "
$ws.Cells.Item(7, 4).Value = "0.636"
$ws.Cells.Item(7, 5).Value = "  +0.07%  "
$ws.Cells.Item(8, 5).Value = "  +0.31%  "
$ws.Cells.Item(9, 4).Value = "0.715"
$ws.Cells.Item(9, 5).Value = "  -0.44%  "
$ws.Cells.Item(10, 4).Value = "0.160"
$ws.Cells.Item(10, 5).Value = "  -3.38%  "
$ws.Cells.Item(11, 4).Value = "56.13"
$ws.Cells.Item(11, 5).Value = "  +6.33%  "
$ws.Cells.Item(12, 4).Value = "0.0000289"
$ws.Cells.Item(12, 5).Value = "  -4.35%  "
$ws.Cells.Item(13, 4).Value = "10.53"
$ws.Cells.Item(13, 5).Value = "  -1.30%  "
$ws.Cells.Item(14, 4).Value = "4.311.61"
$ws.Cells.Item(14, 5).Value = "  +0.22%  "
$ws.Cells.Item(15, 4).Value = "3.719.91"
$ws.Cells.Item(15, 5).Value = "  -0.94%  "
$ws.Cells.Item(16, 4).Value = "19.27"
$ws.Cells.Item(16, 5).Value = "  -0.89%  "
$ws.Cells.Item(17, 2).Value = "TRON"
$ws.Cells.Item(17, 3).Value = "https://coinranking.com/coin/qUhEFk1I61atv+tron-trx"
$ws.Cells.Item(17, 4).Value = "0.126"
$ws.Cells.Item(17, 5).Value = "  -0.56%  "
$ws.Cells.Item(18, 2).Value = "Uniswap"
$ws.Cells.Item(18, 3).Value = "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
$ws.Cells.Item(18, 4).Value = "12.92"
$ws.Cells.Item(18, 5).Value = "  -0.70%  "
$ws.Cells.Item(19, 4).Value = "1.13"
$ws.Cells.Item(19, 5).Value = "  -1.24%  "
$ws.Cells.Item(20, 4).Value = "68.654.26"
$ws.Cells.Item(20, 5).Value = "  +0.97%  "
$ws.Cells.Item(21, 4).Value = "409.82"
$ws.Cells.Item(21, 5).Value = "  -0.08%  "
$ws.Cells.Item(22, 4).Value = "4.60"
$ws.Cells.Item(22, 5).Value = "  -0.18%  "
$ws.Cells.Item(23, 4).Value = "89.15"
$ws.Cells.Item(23, 5).Value = "  +0.50%  "
$ws.Cells.Item(24, 4).Value = "3.02"
$ws.Cells.Item(24, 5).Value = "  -2.61%  "
$ws.Cells.Item(25, 2).Value = "InternetComputer(DFINITY)"
$ws.Cells.Item(25, 3).Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Cells.Item(25, 4).Value = "12.80"
$ws.Cells.Item(25, 5).Value = "  -0.72%  "
$ws.Cells.Item(26, 2).Value = "RenderToken"
$ws.Cells.Item(26, 3).Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Cells.Item(26, 4).Value = "10.98"
$ws.Cells.Item(26, 5).Value = "  +2.07%  "
$ws.Cells.Item(27, 5).Value = "  +1.25%  "
$ws.Cells.Item(28, 4).Value = "3.74"
$ws.Cells.Item(28, 5).Value = "  -3.05%  "
$ws.Cells.Item(29, 4).Value = "9.62"
$ws.Cells.Item(29, 5).Value = "  +0.53%  "
$ws.Cells.Item(30, 4).Value = "33.05"
$ws.Cells.Item(30, 5).Value = "  -0.02%  "
$ws.Cells.Item(31, 4).Value = "7.24"
$ws.Cells.Item(31, 5).Value = "  -11.14%  "
$ws.Cells.Item(32, 4).Value = "12.57"
$ws.Cells.Item(32, 5).Value = "  -1.61%  "
$ws.Cells.Item(33, 4).Value = "0.122"
$ws.Cells.Item(33, 5).Value = "  +3.49%  "
$ws.Cells.Item(34, 4).Value = "624.86"
$ws.Cells.Item(34, 5).Value = "  +5.22%  "
$ws.Cells.Item(35, 4).Value = "44.38"
$ws.Cells.Item(35, 5).Value = "  -0.63%  "
$ws.Cells.Item(36, 4).Value = "65.41"
$ws.Cells.Item(36, 5).Value = "  -1.12%  "
$ws.Cells.Item(37, 2).Value = "PEPE"
$ws.Cells.Item(37, 3).Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Cells.Item(37, 4).Value = "0.0₃0830"
$ws.Cells.Item(37, 5).Value = "  -10.38%  "
$ws.Cells.Item(38, 2).Value = "TheGraph"
$ws.Cells.Item(38, 3).Value = "https://coinranking.com/coin/qhd1biQ7M+thegraph-grt"
$ws.Cells.Item(38, 4).Value = "0.413"
$ws.Cells.Item(38, 5).Value = "  +2.11%  "
$ws.Cells.Item(39, 4).Value = "1.00"
$ws.Cells.Item(39, 5).Value = "  -0.13%  "
$ws.Cells.Item(40, 4).Value = "1.00"
$ws.Cells.Item(40, 5).Value = "  +0.09%  "
$ws.Cells.Item(41, 4).Value = "0.140"
$ws.Cells.Item(41, 5).Value = "  +2.82%  "
$ws.Cells.Item(42, 4).Value = "3.03"
$ws.Cells.Item(42, 5).Value = "  -2.10%  "
$ws.Cells.Item(43, 4).Value = "0.0443"
$ws.Cells.Item(43, 5).Value = "  +0.13%  "
$ws.Cells.Item(44, 4).Value = "2.61"
$ws.Cells.Item(44, 5).Value = "  +0.71%  "
$ws.Cells.Item(45, 4).Value = "0.139"
$ws.Cells.Item(45, 5).Value = "  +2.77%  "
$ws.Cells.Item(46, 4).Value = "2.862.17"
$ws.Cells.Item(46, 5).Value = "  +3.96%  "
$ws.Cells.Item(47, 4).Value = "2.74"
$ws.Cells.Item(47, 5).Value = "  +2.40%  "
$ws.Cells.Item(48, 4).Value = "9.08"
$ws.Cells.Item(48, 5).Value = "  -4.63%  "
$ws.Cells.Item(49, 4).Value = "3.10"
$ws.Cells.Item(49, 5).Value = "  -1.93%  "
$ws.Cells.Item(50, 2).Value = "Monero"
$ws.Cells.Item(50, 3).Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Cells.Item(50, 4).Value = "141.76"
$ws.Cells.Item(50, 5).Value = "  -1.20%  "
$ws.Cells.Item(51, 2).Value = "dogwifhat"
$ws.Cells.Item(51, 3).Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Cells.Item(51, 4).Value = "2.59"
$ws.Cells.Item(51, 5).Value = "  -22.88%  "
